$wb = $excel.ActiveWorkbook

# Sheet "平衡性检验" (Balance test) - update post-match bias / t-stat / p-value values
$ws1 = $wb.Worksheets.Item("平衡性检验")

$ws1.Range("C2").Value = -5.062293621080931
$ws1.Range("D2").Value = 113.1898977236058
$ws1.Range("G2").Value = -1.384053665424951
$ws1.Range("H2").Value = 0.1664453888926849

$ws1.Range("C3").Value = -1.676422489111038
$ws1.Range("D3").Value = 104.6579459240546
$ws1.Range("G3").Value = -0.4583413892060091
$ws1.Range("H3").Value = 0.6467405249150434

$ws1.Range("C4").Value = 16.21594552480293
$ws1.Range("D4").Value = 64.12677777939534
$ws1.Range("G4").Value = 4.433511866730184
$ws1.Range("H4").Value = 0.00000960515358271766

$ws1.Range("C5").Value = 15.99776471493287
$ws1.Range("D5").Value = 64.92584297378158
$ws1.Range("G5").Value = 4.373860259725699
$ws1.Range("H5").Value = 0.00001262426739122109

$ws1.Range("C6").Value = -6.273725120801471
$ws1.Range("D6").Value = 119.2074657667888
$ws1.Range("G6").Value = -1.715264443207027
$ws1.Range("H6").Value = 0.08640043158271631

$ws1.Range("C7").Value = 0.3968688470246407
$ws1.Range("D7").Value = -101.7303809120134
$ws1.Range("G7").Value = 0.1085057137203632
$ws1.Range("H7").Value = 0.9136018339016092

# Sheet "匹配概况" (Matching overview) - restore caliper from 0.02 to 0.05 results
$ws3 = $wb.Worksheets.Item("匹配概况")
$ws3.Range("B7").Value = 2990
